# LOM3246.docx - split the two long "Programa" paragraphs (PT and EN/italic)
# into four sentence groups separated by manual line breaks (w:br), instead
# of being one continuous run of text. No wording changes - only w:br
# elements are introduced at the three sentence-group boundaries of each
# paragraph.

$d = $word.ActiveDocument

# --- Portuguese paragraph ---
# "...picnometria. " | "Análises microestruturais...WDX). " | "Análises térmicas...(TGA)." | "Reometria..."
$d.Content.Find.Execute("picnometria. Análises microestruturais", $true, $false, $false, $false, $false, `
    $true, 1, $false, "picnometria. ^lAnálises microestruturais", 2)

$d.Content.Find.Execute("(EDX e WDX). Análises térmicas", $true, $false, $false, $false, $false, `
    $true, 1, $false, "(EDX e WDX). ^lAnálises térmicas", 2)

$d.Content.Find.Execute("(TGA).Reometria", $true, $false, $false, $false, $false, `
    $true, 1, $false, "(TGA).^lReometria", 2)

# --- English (italic) paragraph ---
# "...pycnometry." | "Microstructural analysis...WDX)." | "Thermal analysis...(TGA)." | "Rheometry..."
$d.Content.Find.Execute("pycnometry.Microstructural analysis", $true, $false, $false, $false, $false, `
    $true, 1, $false, "pycnometry.^lMicrostructural analysis", 2)

$d.Content.Find.Execute("(EDX and WDX).Thermal analysis", $true, $false, $false, $false, $false, `
    $true, 1, $false, "(EDX and WDX).^lThermal analysis", 2)

$d.Content.Find.Execute("(TGA).Rheometry", $true, $false, $false, $false, $false, `
    $true, 1, $false, "(TGA).^lRheometry", 2)
